$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'62.406.56"
$ws.Cells.Item(2, 5).Value = "'  +2.69%  "
$ws.Cells.Item(3, 4).Value = "'2.427.16"
$ws.Cells.Item(3, 5).Value = "'  +3.58%  "
$ws.Cells.Item(4, 5).Value = "'  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'556.89"
$ws.Cells.Item(5, 5).Value = "'  +2.49%  "
$ws.Cells.Item(6, 4).Value = "'143.82"
$ws.Cells.Item(6, 5).Value = "'  +5.49%  "
$ws.Cells.Item(8, 4).Value = "'0.534"
$ws.Cells.Item(8, 5).Value = "'  +1.83%  "
$ws.Cells.Item(9, 4).Value = "'2.427.83"
$ws.Cells.Item(9, 5).Value = "'  +3.65%  "
$ws.Cells.Item(10, 5).Value = "'  +5.41%  "
$ws.Cells.Item(11, 5).Value = "'  -0.26%  "
$ws.Cells.Item(12, 4).Value = "'5.41"
$ws.Cells.Item(12, 5).Value = "'  +2.29%  "
$ws.Cells.Item(13, 5).Value = "'  +2.92%  "
$ws.Cells.Item(14, 4).Value = "'26.36"
$ws.Cells.Item(14, 5).Value = "'  +7.26%  "
$ws.Cells.Item(15, 4).Value = "'0.0000175"
$ws.Cells.Item(15, 5).Value = "'  +9.85%  "
$ws.Cells.Item(16, 4).Value = "'2.864.33"
$ws.Cells.Item(16, 5).Value = "'  +3.62%  "
$ws.Cells.Item(17, 4).Value = "'62.206.44"
$ws.Cells.Item(17, 5).Value = "'  +2.53%  "
$ws.Cells.Item(18, 4).Value = "'2.425.75"
$ws.Cells.Item(18, 5).Value = "'  +3.37%  "
$ws.Cells.Item(19, 4).Value = "'11.13"
$ws.Cells.Item(19, 5).Value = "'  +4.83%  "
$ws.Cells.Item(20, 4).Value = "'324.80"
$ws.Cells.Item(20, 5).Value = "'  +2.00%  "
$ws.Cells.Item(21, 4).Value = "'4.18"
$ws.Cells.Item(21, 5).Value = "'  +1.52%  "
$ws.Cells.Item(22, 5).Value = "'  +3.58%  "
$ws.Cells.Item(23, 5).Value = "'  +0.31%  "
$ws.Cells.Item(24, 5).Value = "'  +5.69%  "
$ws.Cells.Item(25, 4).Value = "'65.05"
$ws.Cells.Item(25, 5).Value = "'  +2.88%  "
$ws.Cells.Item(26, 4).Value = "'9.10"
$ws.Cells.Item(26, 5).Value = "'  +7.24%  "
$ws.Cells.Item(27, 4).Value = "'569.28"
$ws.Cells.Item(27, 5).Value = "'  +14.47%  "
$ws.Cells.Item(28, 4).Value = "'2.544.83"
$ws.Cells.Item(28, 5).Value = "'  +3.53%  "
$ws.Cells.Item(29, 4).Value = "'0.999"
$ws.Cells.Item(29, 5).Value = "'  -0.23%  "
$ws.Cells.Item(30, 4).Value = "'0.0₃0946"
$ws.Cells.Item(30, 5).Value = "'  +10.14%  "
$ws.Cells.Item(31, 4).Value = "'8.42"
$ws.Cells.Item(31, 5).Value = "'  +6.21%  "
$ws.Cells.Item(32, 5).Value = "'  +6.17%  "
$ws.Cells.Item(33, 5).Value = "'  +2.31%  "
$ws.Cells.Item(34, 5).Value = "'  +4.36%  "
$ws.Cells.Item(35, 5).Value = "'  +5.37%  "
$ws.Cells.Item(36, 5).Value = "'  +9.45%  "
$ws.Cells.Item(37, 5).Value = "'  +5.85%  "
$ws.Cells.Item(39, 4).Value = "'0.386"
$ws.Cells.Item(39, 5).Value = "'  +2.87%  "
$ws.Cells.Item(40, 5).Value = "'  +3.70%  "
$ws.Cells.Item(41, 4).Value = "'18.82"
$ws.Cells.Item(41, 5).Value = "'  +1.89%  "
$ws.Cells.Item(42, 4).Value = "'149.91"
$ws.Cells.Item(42, 5).Value = "'  +5.08%  "
$ws.Cells.Item(43, 5).Value = "'  +0.02%  "
$ws.Cells.Item(44, 4).Value = "'41.70"
$ws.Cells.Item(44, 5).Value = "'  +2.95%  "
$ws.Cells.Item(45, 4).Value = "'2.35"
$ws.Cells.Item(45, 5).Value = "'  +15.61%  "
$ws.Cells.Item(46, 4).Value = "'151.48"
$ws.Cells.Item(46, 5).Value = "'  +6.41%  "
$ws.Cells.Item(47, 4).Value = "'3.65"
$ws.Cells.Item(47, 5).Value = "'  +2.89%  "
$ws.Cells.Item(48, 5).Value = "'  +5.10%  "
$ws.Cells.Item(49, 4).Value = "'20.45"
$ws.Cells.Item(49, 5).Value = "'  +7.62%  "
$ws.Cells.Item(50, 5).Value = "'  +4.35%  "
$ws.Cells.Item(51, 5).Value = "'  +4.00%  "
